# Generated edit script: apply row content swaps + append new row 122
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update swapped match rows (F:V) ---
# Row 11
$ws.Range("F11").Value = "Selimbar"
$ws.Range("G11").Value = 2
$ws.Range("H11").Value = "Unirea Dej"
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 2.06
$ws.Range("K11").Value = "10/08/2023 22:12"
$ws.Range("L11").Value = 2.2
$ws.Range("M11").Value = "12/08/2023 09:56"
$ws.Range("N11").Value = 3
$ws.Range("O11").Value = "10/08/2023 22:12"
$ws.Range("P11").Value = 2.98
$ws.Range("Q11").Value = "12/08/2023 09:56"
$ws.Range("R11").Value = 3.52
$ws.Range("S11").Value = "10/08/2023 22:12"
$ws.Range("T11").Value = 3.73
$ws.Range("U11").Value = "12/08/2023 09:56"
$ws.Range("V11").Value = "https://www.betexplorer.com/football/romania/liga-2/selimbar-unirea-dej/OYyYABgA/"

# Row 16
$ws.Range("F16").Value = "Unirea Slobozia"
$ws.Range("G16").Value = 1
$ws.Range("H16").Value = "Viitorul Tg. Jiu"
$ws.Range("I16").Value = 1
$ws.Range("J16").Value = 1.58
$ws.Range("K16").Value = "10/08/2023 22:12"
$ws.Range("L16").Value = 1.49
$ws.Range("M16").Value = "12/08/2023 09:14"
$ws.Range("N16").Value = 3.77
$ws.Range("O16").Value = "10/08/2023 22:12"
$ws.Range("P16").Value = 4.35
$ws.Range("Q16").Value = "12/08/2023 09:14"
$ws.Range("R16").Value = 5.17
$ws.Range("S16").Value = "10/08/2023 22:12"
$ws.Range("T16").Value = 6.51
$ws.Range("U16").Value = "12/08/2023 09:14"
$ws.Range("V16").Value = "https://www.betexplorer.com/football/romania/liga-2/unirea-slobozia-viitorul-targu-jiu/rwLG2m1j/"

# Row 83
$ws.Range("F83").Value = "Viitorul Tg. Jiu"
$ws.Range("G83").Value = 2
$ws.Range("H83").Value = "Concordia"
$ws.Range("I83").Value = 2
$ws.Range("J83").Value = 4.16
$ws.Range("K83").Value = "05/10/2023 21:13"
$ws.Range("L83").Value = 5.13
$ws.Range("M83").Value = "07/10/2023 09:52"
$ws.Range("N83").Value = 3.41
$ws.Range("O83").Value = "05/10/2023 21:13"
$ws.Range("P83").Value = 3.55
$ws.Range("Q83").Value = "07/10/2023 09:52"
$ws.Range("R83").Value = 1.76
$ws.Range("S83").Value = "05/10/2023 21:13"
$ws.Range("T83").Value = 1.72
$ws.Range("U83").Value = "07/10/2023 09:52"
$ws.Range("V83").Value = "https://www.betexplorer.com/football/romania/liga-2/viitorul-targu-jiu-concordia/21fYhhjF/"

# Row 84
$ws.Range("F84").Value = "Unirea Dej"
$ws.Range("G84").Value = 0
$ws.Range("H84").Value = "CSC Dumbravita"
$ws.Range("I84").Value = 1
$ws.Range("J84").Value = 2.17
$ws.Range("K84").Value = "05/10/2023 21:13"
$ws.Range("L84").Value = 2.41
$ws.Range("M84").Value = "07/10/2023 09:44"
$ws.Range("N84").Value = 3.2
$ws.Range("O84").Value = "05/10/2023 21:13"
$ws.Range("P84").Value = 3.14
$ws.Range("Q84").Value = "07/10/2023 09:43"
$ws.Range("R84").Value = 3.04
$ws.Range("S84").Value = "05/10/2023 21:13"
$ws.Range("T84").Value = 3.08
$ws.Range("U84").Value = "07/10/2023 09:44"
$ws.Range("V84").Value = "https://www.betexplorer.com/football/romania/liga-2/unirea-dej-csc-dumbravita/hphGpUce/"

# Row 85
$ws.Range("F85").Value = "Mioveni"
$ws.Range("G85").Value = 2
$ws.Range("H85").Value = "CSM Resita"
$ws.Range("I85").Value = 1
$ws.Range("J85").Value = 1.75
$ws.Range("K85").Value = "05/10/2023 21:13"
$ws.Range("L85").Value = 1.81
$ws.Range("M85").Value = "07/10/2023 09:51"
$ws.Range("N85").Value = 3.38
$ws.Range("O85").Value = "05/10/2023 21:13"
$ws.Range("P85").Value = 3.47
$ws.Range("Q85").Value = "07/10/2023 09:51"
$ws.Range("R85").Value = 4.26
$ws.Range("S85").Value = "05/10/2023 21:13"
$ws.Range("T85").Value = 4.59
$ws.Range("U85").Value = "07/10/2023 09:51"
$ws.Range("V85").Value = "https://www.betexplorer.com/football/romania/liga-2/mioveni-csm-resita/thgxhC5L/"

# Row 86
$ws.Range("F86").Value = "Alexandria"
$ws.Range("G86").Value = 2
$ws.Range("H86").Value = "Ceahlaul"
$ws.Range("I86").Value = 2
$ws.Range("J86").Value = 3.31
$ws.Range("K86").Value = "05/10/2023 21:13"
$ws.Range("L86").Value = 3.37
$ws.Range("M86").Value = "07/10/2023 09:53"
$ws.Range("N86").Value = 3.06
$ws.Range("O86").Value = "05/10/2023 21:13"
$ws.Range("P86").Value = 3.19
$ws.Range("Q86").Value = "07/10/2023 09:36"
$ws.Range("R86").Value = 2.12
$ws.Range("S86").Value = "05/10/2023 21:13"
$ws.Range("T86").Value = 2.24
$ws.Range("U86").Value = "07/10/2023 09:53"
$ws.Range("V86").Value = "https://www.betexplorer.com/football/romania/liga-2/csm-alexandria-ceahlaul/z9sLqlD1/"

# Row 87
$ws.Range("F87").Value = "Chindia Targoviste"
$ws.Range("G87").Value = 3
$ws.Range("H87").Value = "Progresul Spartac"
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 1.2
$ws.Range("K87").Value = "07/10/2023 00:13"
$ws.Range("L87").Value = 1.17
$ws.Range("M87").Value = "07/10/2023 08:28"
$ws.Range("N87").Value = 6.4
$ws.Range("O87").Value = "07/10/2023 00:13"
$ws.Range("P87").Value = 7.19
$ws.Range("Q87").Value = "07/10/2023 09:33"
$ws.Range("R87").Value = 12.64
$ws.Range("S87").Value = "07/10/2023 00:13"
$ws.Range("T87").Value = 16.45
$ws.Range("U87").Value = "07/10/2023 09:33"
$ws.Range("V87").Value = "https://www.betexplorer.com/football/romania/liga-2/chindia-targoviste-progresul-spartac/0xz9njTr/"

# Row 88
$ws.Range("F88").Value = "Metaloglobus Bucharest"
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = "Csikszereda M. Ciuc"
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 2.48
$ws.Range("K88").Value = "05/10/2023 21:13"
$ws.Range("L88").Value = 2.46
$ws.Range("M88").Value = "07/10/2023 09:57"
$ws.Range("N88").Value = 2.95
$ws.Range("O88").Value = "05/10/2023 21:13"
$ws.Range("P88").Value = 3.13
$ws.Range("Q88").Value = "07/10/2023 09:57"
$ws.Range("R88").Value = 2.79
$ws.Range("S88").Value = "05/10/2023 21:13"
$ws.Range("T88").Value = 3.03
$ws.Range("U88").Value = "07/10/2023 09:57"
$ws.Range("V88").Value = "https://www.betexplorer.com/football/romania/liga-2/metaloglobus-bucharest-miercurea-ciuc/faZCoArk/"

# Row 94
$ws.Range("F94").Value = "CSC Dumbravita"
$ws.Range("G94").Value = 2
$ws.Range("H94").Value = "Alexandria"
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 2.07
$ws.Range("K94").Value = "20/10/2023 04:42"
$ws.Range("L94").Value = 2.05
$ws.Range("M94").Value = "21/10/2023 09:51"
$ws.Range("N94").Value = 3.2
$ws.Range("O94").Value = "20/10/2023 04:42"
$ws.Range("P94").Value = 3.41
$ws.Range("Q94").Value = "21/10/2023 09:51"
$ws.Range("R94").Value = 3.38
$ws.Range("S94").Value = "20/10/2023 04:42"
$ws.Range("T94").Value = 3.63
$ws.Range("U94").Value = "21/10/2023 09:51"
$ws.Range("V94").Value = "https://www.betexplorer.com/football/romania/liga-2/csc-dumbravita-csm-alexandria/EFUnbmSf/"

# Row 95
$ws.Range("F95").Value = "CSM Resita"
$ws.Range("G95").Value = 1
$ws.Range("H95").Value = "Tunari"
$ws.Range("I95").Value = 1
$ws.Range("J95").Value = 1.92
$ws.Range("K95").Value = "19/10/2023 21:12"
$ws.Range("L95").Value = 1.85
$ws.Range("M95").Value = "21/10/2023 06:21"
$ws.Range("N95").Value = 3.35
$ws.Range("O95").Value = "19/10/2023 21:12"
$ws.Range("P95").Value = 3.77
$ws.Range("Q95").Value = "21/10/2023 08:01"
$ws.Range("R95").Value = 3.54
$ws.Range("S95").Value = "19/10/2023 21:12"
$ws.Range("T95").Value = 3.9
$ws.Range("U95").Value = "21/10/2023 06:21"
$ws.Range("V95").Value = "https://www.betexplorer.com/football/romania/liga-2/csm-resita-tunari/2yBsKUlE/"

# Row 96
$ws.Range("F96").Value = "Progresul Spartac"
$ws.Range("G96").Value = 0
$ws.Range("H96").Value = "Metaloglobus Bucharest"
$ws.Range("I96").Value = 1
$ws.Range("J96").Value = 3.91
$ws.Range("K96").Value = "19/10/2023 21:12"
$ws.Range("L96").Value = 5.54
$ws.Range("M96").Value = "21/10/2023 09:58"
$ws.Range("N96").Value = 3.34
$ws.Range("O96").Value = "19/10/2023 21:12"
$ws.Range("P96").Value = 3.52
$ws.Range("Q96").Value = "21/10/2023 09:58"
$ws.Range("R96").Value = 1.83
$ws.Range("S96").Value = "19/10/2023 21:12"
$ws.Range("T96").Value = 1.68
$ws.Range("U96").Value = "21/10/2023 09:58"
$ws.Range("V96").Value = "https://www.betexplorer.com/football/romania/liga-2/progresul-spartac-metaloglobus-bucharest/nTSfdRd7/"

# Row 104
$ws.Range("F104").Value = "Selimbar"
$ws.Range("G104").Value = 1
$ws.Range("H104").Value = "CSC Dumbravita"
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 1.75
$ws.Range("K104").Value = "26/10/2023 22:12"
$ws.Range("L104").Value = 1.79
$ws.Range("M104").Value = "28/10/2023 09:51"
$ws.Range("N104").Value = 3.4
$ws.Range("O104").Value = "26/10/2023 22:12"
$ws.Range("P104").Value = 3.56
$ws.Range("Q104").Value = "28/10/2023 09:58"
$ws.Range("R104").Value = 4.23
$ws.Range("S104").Value = "26/10/2023 22:12"
$ws.Range("T104").Value = 4.62
$ws.Range("U104").Value = "28/10/2023 09:51"
$ws.Range("V104").Value = "https://www.betexplorer.com/football/romania/liga-2/selimbar-csc-dumbravita/dbnw94tD/"

# Row 105
$ws.Range("F105").Value = "Mioveni"
$ws.Range("G105").Value = 1
$ws.Range("H105").Value = "Concordia"
$ws.Range("I105").Value = 1
$ws.Range("J105").Value = 3.02
$ws.Range("K105").Value = "26/10/2023 22:12"
$ws.Range("L105").Value = 3.13
$ws.Range("M105").Value = "28/10/2023 09:52"
$ws.Range("N105").Value = 3.04
$ws.Range("O105").Value = "26/10/2023 22:12"
$ws.Range("P105").Value = 2.92
$ws.Range("Q105").Value = "28/10/2023 09:09"
$ws.Range("R105").Value = 2.27
$ws.Range("S105").Value = "26/10/2023 22:12"
$ws.Range("T105").Value = 2.53
$ws.Range("U105").Value = "28/10/2023 09:52"
$ws.Range("V105").Value = "https://www.betexplorer.com/football/romania/liga-2/mioveni-concordia/A1EkI8JQ/"

# Row 107
$ws.Range("F107").Value = "Tunari"
$ws.Range("G107").Value = 0
$ws.Range("H107").Value = "Viitorul Tg. Jiu"
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 2.13
$ws.Range("K107").Value = "26/10/2023 22:12"
$ws.Range("L107").Value = 2.18
$ws.Range("M107").Value = "28/10/2023 09:59"
$ws.Range("N107").Value = 3.22
$ws.Range("O107").Value = "26/10/2023 22:12"
$ws.Range("P107").Value = 3.46
$ws.Range("Q107").Value = "28/10/2023 09:59"
$ws.Range("R107").Value = 3.11
$ws.Range("S107").Value = "26/10/2023 22:12"
$ws.Range("T107").Value = 3.24
$ws.Range("U107").Value = "28/10/2023 09:59"
$ws.Range("V107").Value = "https://www.betexplorer.com/football/romania/liga-2/tunari-viitorul-targu-jiu/CA0HDSRs/"
# --- Append new row 122 ---
$ws.Range("A121").Copy()
$ws.Range("A122").PasteSpecial(-4122)
$ws.Range("E121").Copy()
$ws.Range("E122").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A122").Value = 121
$ws.Range("B122").Value = "romania"
$ws.Range("C122").Value = "liga-2"
$ws.Range("D122").Value = "2023-2024"
$ws.Range("E122").Value = 45239.77083333334
$ws.Range("F122").Value = "Chindia Targoviste"
$ws.Range("G122").Value = 0
$ws.Range("H122").Value = "Mioveni"
$ws.Range("I122").Value = 2
$ws.Range("J122").Value = 2.13
$ws.Range("K122").Value = "09/11/2023 08:15"
$ws.Range("L122").Value = 1.83
$ws.Range("M122").Value = "09/11/2023 18:25"
$ws.Range("N122").Value = 3.02
$ws.Range("O122").Value = "09/11/2023 08:15"
$ws.Range("P122").Value = 3.23
$ws.Range("Q122").Value = "09/11/2023 18:25"
$ws.Range("R122").Value = 3.79
$ws.Range("S122").Value = "09/11/2023 08:15"
$ws.Range("T122").Value = 4.94
$ws.Range("U122").Value = "09/11/2023 18:28"
$ws.Range("V122").Value = "https://www.betexplorer.com/football/romania/liga-2/chindia-targoviste-mioveni/AsGuA5r4/"
